# Updates cryptocurrency price/volume data to reflect the latest scrape
# (cryptos.xlsx, GitHub Actions run on Wed Jul 24 12:22:55 UTC 2024).
#
# Column D ("Price") cells are stored as text in this sheet (values like
# "66.399.81" use '.' as both thousands- and decimal-separator, and values
# such as "1.00" / "0.0291" must keep their exact trailing zeros), so each
# Price cell's NumberFormat is forced to "@" (Text) before the new value is
# assigned -- otherwise Excel would auto-convert the numeric-looking text
# into a floating point number and corrupt the formatting/precision.
# Column E ("Volume(1h)") values already contain padding spaces and a "%"
# sign, so Excel keeps them as plain text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.376.51"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.457.80"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.70"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.10"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  +5.22%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.455.51"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.418"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.061.37"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.03"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.276.31"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.451.38"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.98"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.86"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.09"
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.66"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.536"
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("E26").Value = "  +4.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.03"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("E28").Value = "  +3.13%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.96"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.67"
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.07"
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("E35").Value = "  -5.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.56"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.42"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.884"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.04"
$ws.Range("E39").Value = "  -6.53%  "
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.821.88"
$ws.Range("E41").Value = "  +3.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.52"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.57"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.43"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.25"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.41"
$ws.Range("E47").Value = "  +5.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.08"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0290"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("E50").Value = "  +2.90%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.995"
$ws.Range("E51").Value = "  -2.01%  "

# Row 51 previously listed Arweave; it now lists ONDO with a new link,
# price and volume figure (handled above via B51/C51/D51/E51).
